# Update the validation-report workbook: rename "Obs" -> "Obs_relatorio",
# add a new "Obs_sped" column, and replace the boolean "VERDADEIRO" flags
# with the full success message on both sheets.

$wb = $excel.ActiveWorkbook

$msg = "Validado com sucesso! Nenhuma divergência entre o SPED e o relatório foi encontrada!"

# --- Sheet "Bico" (H = Obs_relatorio, I = new Obs_sped) ---
$ws1 = $wb.Worksheets.Item("Bico")
$ws1.Cells.Item(1, 8).Value = "Obs_relatorio"
$ws1.Cells.Item(1, 9).Value = "Obs_sped"

for ($r = 2; $r -le 13; $r++) {
    $ws1.Cells.Item($r, 8).Value = $msg
    # Create an empty (but present) text cell in the new Obs_sped column,
    # matching the blank-but-existing cells already used elsewhere in the file.
    $ws1.Cells.Item($r, 9).Value = "'"
}

# --- Sheet "Tanque" (F = Obs_relatorio, G = new Obs_sped) ---
$ws2 = $wb.Worksheets.Item("Tanque")
$ws2.Cells.Item(1, 6).Value = "Obs_relatorio"
$ws2.Cells.Item(1, 7).Value = "Obs_sped"

for ($r = 2; $r -le 6; $r++) {
    $ws2.Cells.Item($r, 6).Value = $msg
    $ws2.Cells.Item($r, 7).Value = "'"
}
